# Update vm_pu.xlsx Case_5_133 result values ("case with 380 kV done")
# Rows 2-25 correspond to bus indices 0-23. Columns B-F and I-M hold
# voltage magnitude (p.u.) results that were recomputed; G, H, N are
# unchanged (slack/reference values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M")

$data = @(
    @(1.02,1.067515367707383,1.066648265962144,1.071631246931134,1.081017199031096,1.056387605843958,1.07245987330346,1.06935824989916,1.07432790634683,1.08368912070072),
    @(1.02,1.068928090650741,1.067776926209493,1.072950565568604,1.082508504673063,1.056912375903674,1.073526819397444,1.070301762635523,1.075462574503108,1.084997162449742),
    @(1.02,1.069840582125223,1.068505779158423,1.073802457860655,1.083472070239699,1.057249789916708,1.074215097625614,1.070910227125933,1.076194438545745,1.0858416276603),
    @(1.02,1.070223809267357,1.068811843006532,1.074160170210607,1.083876823210907,1.057391127720013,1.0745039507973,1.071165539234296,1.076501559590203,1.08619618661742),
    @(1.02,1.070288132461976,1.06886321233133,1.074220207063781,1.083944763876282,1.05741482907846,1.0745524214132,1.071208378877077,1.076553094219964,1.086255692184816),
    @(1.02,1.069845704326204,1.068509870151678,1.073807239285192,1.083477479856299,1.057251680483691,1.074218959249759,1.070913640524047,1.076198544483796,1.085846367070903),
    @(1.02,1.067993146799555,1.067030008455863,1.072077492496162,1.081521488884436,1.056565400780958,1.072820892358421,1.069677542046553,1.074711861403927,1.084131582214222),
    @(1.02,1.06471586511701,1.064410854184382,1.069015413720006,1.078063664197841,1.055339506423568,1.070340934330664,1.067483454796774,1.072073923242188,1.081094865979735),
    @(1.02,1.062521931227502,1.062656738991818,1.066964169941332,1.075750511212564,1.054510903166256,1.068676265371207,1.066009723804425,1.070302666590071,1.079059849365235),
    @(1.02,1.061569679915805,1.06189521655306,1.066073527469754,1.07474690480448,1.054149375934798,1.067952672221848,1.06536890413672,1.069532610462937,1.078176070253449),
    @(1.02,1.061215622603977,1.061612049704286,1.065742328453265,1.074373811865063,1.054014673706317,1.067683472977031,1.065130465984984,1.069246105508213,1.077847395950545),
    @(1.02,1.061291585022746,1.061672803788889,1.065813388855501,1.07445385568724,1.054043586628528,1.067741236449419,1.065181630389945,1.069307583266172,1.077917915910562),
    @(1.02,1.061540420593538,1.061871816120465,1.066046158147569,1.074716071177981,1.054138249892751,1.067930428828082,1.065349203146411,1.069508937542846,1.078148910128408),
    @(1.02,1.061693689923255,1.061994393782384,1.06618952508961,1.074877589602975,1.054196519959848,1.068046940038117,1.065452395852965,1.06963293581071,1.078291180118931),
    @(1.02,1.062585080579933,1.062707236506135,1.067023226842336,1.075817074448948,1.05453483855984,1.06872422870491,1.066052195835897,1.070353706927716,1.079118447319322),
    @(1.02,1.063143614263222,1.063153849730672,1.067545526862394,1.076405847799309,1.054746321566114,1.069148324793548,1.066427711678791,1.070804994939156,1.079636667348151),
    @(1.02,1.063469179712746,1.063414161146435,1.067849940717619,1.076749077367921,1.054869412245331,1.069395425050574,1.066646485109741,1.071067925943711,1.079938685685199),
    @(1.02,1.063580152437736,1.063502888482465,1.067953698374791,1.076866077490464,1.054911338358063,1.069479634630706,1.066721037509099,1.071157528360652,1.080041623868967),
    @(1.02,1.063083711475979,1.063105952109838,1.067489513369828,1.076342697911816,1.054723658752446,1.069102851047566,1.066387449174326,1.070756606849202,1.079581093227883),
    @(1.02,1.061467154418257,1.06181322037734,1.065977623802941,1.074638863849232,1.05411038540514,1.067874728162014,1.065299868489596,1.069449656812422,1.078080899180301),
    @(1.02,1.060448738560576,1.060998669520466,1.065024866596093,1.073565804437943,1.053722393629688,1.067100098950237,1.064613692920313,1.068625191857733,1.077135352759878),
    @(1.02,1.060988814690828,1.061430647163649,1.065530149842573,1.074134826342886,1.053928304408475,1.067510979975345,1.06497767424015,1.069062517973675,1.077636827000119),
    @(1.02,1.063110779640251,1.063127595571906,1.067514824182121,1.076371233221114,1.054733899917605,1.069123399503704,1.066405642866278,1.070778472264133,1.079606205555335),
    @(1.02,1.065564688301233,1.065089357342456,1.069808741361133,1.078958959640562,1.055658414171002,1.070984040247775,1.068052596723487,1.072758091355677,1.081881756323601)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $vals[$j]
    }
}
